$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.9931674706421231
$ws.Range("D2").Value = 0.9911988362902867
$ws.Range("E2").Value = 0.9899038524349963
$ws.Range("F2").Value = 0.989083283285945
$ws.Range("G2").Value = 0.9888296519732972
$ws.Range("H2").Value = 0.9905303689921979
$ws.Range("I2").Value = 0.989432384460235
$ws.Range("J2").Value = 0.9885176257931687
$ws.Range("K2").Value = 0.9885765155212325

$ws.Range("C3").Value = 0.9934611064312664
$ws.Range("D3").Value = 0.9919198326531968
$ws.Range("E3").Value = 0.991052411202241
$ws.Range("F3").Value = 0.9904409546753616
$ws.Range("G3").Value = 0.9901876735579727
$ws.Range("H3").Value = 0.9912518525686499
$ws.Range("I3").Value = 0.990799843963233
$ws.Range("J3").Value = 0.9900842022715296
$ws.Range("K3").Value = 0.9899348851777835
